# Updated cryptos list on Sun Oct 20 09:48:17 UTC 2024 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h)/% change (column E) figures
# for the cryptocurrency table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new Price (column D) text. $null means the Price cell is unchanged.
$dUpdates = @{
    2  = "68.428.57"
    3  = "2.648.74"
    4  = "1.00"
    5  = "597.19"
    6  = "159.26"
    8  = "0.539"
    9  = "2.648.42"
    14 = "28.00"
    15 = "3.131.49"
    17 = "68.302.24"
    18 = "2.664.35"
    19 = "11.47"
    20 = "363.88"
    21 = "7.44"
    23 = "4.76"
    25 = "74.57"
    27 = "9.83"
    31 = "562.24"
    39 = "19.64"
    40 = "0.370"
    42 = "5.33"
    46 = "158.15"
    50 = "0.0778"
    51 = "0.575"
}

# Row -> new Volume(1h) (column E) text (percentages keep their padding spaces).
$eUpdates = @{
    2  = "  +0.09%  "
    3  = "  +0.10%  "
    4  = "  +0.07%  "
    5  = "  -0.17%  "
    6  = "  +2.90%  "
    7  = "  -0.01%  "
    8  = "  -1.41%  "
    9  = "  +0.14%  "
    10 = "  -1.56%  "
    11 = "  -1.13%  "
    12 = "  +0.46%  "
    13 = "  -1.22%  "
    14 = "  -0.13%  "
    15 = "  +0.11%  "
    16 = "  -2.98%  "
    17 = "  +0.05%  "
    18 = "  +0.82%  "
    19 = "  +0.93%  "
    20 = "  -0.09%  "
    21 = "  -0.71%  "
    22 = "  +0.71%  "
    23 = "  -2.53%  "
    24 = "  +0.26%  "
    25 = "  -0.36%  "
    26 = "  +0.10%  "
    27 = "  +0.13%  "
    29 = "  -2.96%  "
    30 = "  +0.00%  "
    31 = "  -1.73%  "
    32 = "  -0.38%  "
    33 = "  -1.49%  "
    34 = "  -0.26%  "
    35 = "  +3.70%  "
    36 = "  -1.37%  "
    37 = "  +0.00%  "
    38 = "  -0.47%  "
    39 = "  +1.42%  "
    40 = "  -1.24%  "
    41 = "  -1.15%  "
    42 = "  -0.84%  "
    43 = "  -0.91%  "
    44 = "  -5.08%  "
    45 = "  +0.06%  "
    46 = "  +1.01%  "
    47 = "  +1.89%  "
    48 = "  +0.19%  "
    49 = "  -1.06%  "
    50 = "  -1.30%  "
    51 = "  +1.54%  "
}

# Some Price figures (e.g. "1.00", "597.19") look like plain numbers, and a
# bare Range.Value assignment would let Excel auto-convert them into numeric
# values (losing the fixed-decimal text form the sheet actually stores).
# Force those specific cells to Text format just long enough to type the
# value in as a string, then restore the Normal style so no stray
# number-format/style change is left behind on the cell.
foreach ($row in $dUpdates.Keys) {
    $text = $dUpdates[$row]
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

foreach ($row in $eUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $eUpdates[$row]
}
